# Reseller add & refactoring old code
#
# - createOrgTest1!E2: "AgreementName-1" -> "autoTestPayment1"
# - Make the "createOrgTest1" sheet the active tab / active sheet
# - Update the "createOrgTest1" sheet's selection to F2 (was H6)
# - "TestCases" sheet loses its tabSelected flag (handled automatically
#   once a different sheet becomes active)

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("createOrgTest1")

# Update the payment agreement name value.
$ws2.Range("E2").Value = "autoTestPayment1"

# Switch the active sheet to createOrgTest1 and move the selection to F2.
$ws2.Activate()
$ws2.Range("F2").Select() | Out-Null
